$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) to remain text so that values such as
# "1.017" or "334.50" are not reinterpreted as numbers (which would
# drop trailing zeros / change precision).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.840.86"
$ws.Range("E2").Value = "  +1.35%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.887.75"
$ws.Range("E3").Value = "  +1.23%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.017"
$ws.Range("E4").Value = "  +1.52%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "334.50"
$ws.Range("E5").Value = "  +1.51%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.016"
$ws.Range("E6").Value = "  +1.49%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4686"
$ws.Range("E7").Value = "  -0.58%  "
$ws.Range("E8").Value = "  -1.57%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.65"
$ws.Range("E9").Value = "  +0.44%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08064"
$ws.Range("E10").Value = "  +0.50%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.018"
$ws.Range("E11").Value = "  -0.16%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.96"
$ws.Range("E12").Value = "  +1.50%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.910.20"
$ws.Range("E13").Value = "  +2.89%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.955"
$ws.Range("E14").Value = "  -0.01%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.086"
$ws.Range("E15").Value = "  -1.51%  "
$ws.Range("E16").Value = "  +1.45%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.06766"
$ws.Range("E17").Value = "  +2.91%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "87.29"
$ws.Range("E18").Value = "  +0.63%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.00001048"
$ws.Range("E19").Value = "  +0.78%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.20"
$ws.Range("E20").Value = "  -0.52%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.016"
$ws.Range("E21").Value = "  +1.51%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "27.869.70"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.503"
$ws.Range("E23").Value = "  -0.21%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.94"
$ws.Range("E24").Value = "  -0.29%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.345"
$ws.Range("E25").Value = "  +1.92%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.126.26"
$ws.Range("E26").Value = "  +2.32%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "160.13"
$ws.Range("E27").Value = "  +3.99%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.04"
$ws.Range("E28").Value = "  -1.09%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.090"
$ws.Range("E29").Value = "  +0.40%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.505"
$ws.Range("E30").Value = "  -0.76%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "121.88"
$ws.Range("E31").Value = "  -0.39%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9724"
$ws.Range("E32").Value = "  +1.62%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09468"
$ws.Range("E33").Value = "  -0.25%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.644"
$ws.Range("E34").Value = "  +1.36%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.414"
$ws.Range("E35").Value = "  -3.66%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.353"
$ws.Range("E36").Value = "  +1.05%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06126"
$ws.Range("E37").Value = "  +0.41%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02260"
$ws.Range("E38").Value = "  +0.51%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.215"
$ws.Range("E39").Value = "  -0.39%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5987"
$ws.Range("E40").Value = "  +0.10%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.012"
$ws.Range("E41").Value = "  -1.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1889"
$ws.Range("E42").Value = "  -0.42%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.27"
$ws.Range("E43").Value = "  -0.54%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.265"
$ws.Range("E44").Value = "  -0.09%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5683"
$ws.Range("E45").Value = "  +0.14%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.15"
$ws.Range("E46").Value = "  -0.48%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.410"
$ws.Range("E47").Value = "  -0.18%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.933"
$ws.Range("E48").Value = "  -0.65%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06908"
$ws.Range("E49").Value = "  +1.89%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "113.99"
$ws.Range("E50").Value = "  +3.51%  "
$ws.Range("E51").Value = "  +0.51%  "
